# Update cryptos list with refreshed prices / volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 / 28 swap: PEPE <-> Bittensor (now row 27 = PEPE, row 28 = Bittensor) ---
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

# --- Volume(1h) column (E): plain text percentages, safe to assign directly ---
$ws.Range('E2').Value = '  -3.92%  '
$ws.Range('E3').Value = '  -6.54%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -5.35%  '
$ws.Range('E6').Value = '  -6.63%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -7.09%  '
$ws.Range('E9').Value = '  -6.57%  '
$ws.Range('E10').Value = '  -10.15%  '
$ws.Range('E11').Value = '  -6.38%  '
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('E13').Value = '  -8.55%  '
$ws.Range('E14').Value = '  -9.12%  '
$ws.Range('E15').Value = '  -6.50%  '
$ws.Range('E16').Value = '  -9.58%  '
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('E18').Value = '  -6.87%  '
$ws.Range('E19').Value = '  -8.75%  '
$ws.Range('E20').Value = '  -8.25%  '
$ws.Range('E21').Value = '  -7.96%  '
$ws.Range('E22').Value = '  -7.16%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -2.74%  '
$ws.Range('E26').Value = '  -6.80%  '
$ws.Range('E27').Value = '  -14.35%  '
$ws.Range('E28').Value = '  -6.16%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -11.44%  '
$ws.Range('E31').Value = '  -10.65%  '
$ws.Range('E32').Value = '  -7.96%  '
$ws.Range('E33').Value = '  -8.31%  '
$ws.Range('E34').Value = '  -7.67%  '
$ws.Range('E35').Value = '  -8.69%  '
$ws.Range('E36').Value = '  -11.90%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  -11.83%  '
$ws.Range('E39').Value = '  -6.36%  '
$ws.Range('E40').Value = '  -6.94%  '
$ws.Range('E41').Value = '  -7.58%  '
$ws.Range('E42').Value = '  -8.28%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('E45').Value = '  -10.40%  '
$ws.Range('E46').Value = '  -10.21%  '
$ws.Range('E47').Value = '  -8.50%  '
$ws.Range('E48').Value = '  -11.20%  '
$ws.Range('E50').Value = '  -7.44%  '
$ws.Range('E51').Value = '  -6.65%  '

# --- Price column (D): values look numeric, so force text storage the way Excel
#     does for quote-prefixed numeric text (NumberFormat "@"), then drop the
#     leftover number-format back to Normal so no visible style change remains. ---
$dPrices = [ordered]@{
  'D2' = '61.280.17'
  'D3' = '2.458.37'
  'D5' = '547.39'
  'D6' = '146.36'
  'D8' = '0.586'
  'D9' = '2.456.74'
  'D11' = '5.44'
  'D14' = '26.05'
  'D15' = '2.904.13'
  'D17' = '61.214.73'
  'D18' = '2.469.70'
  'D19' = '11.10'
  'D21' = '4.16'
  'D22' = '318.88'
  'D25' = '63.85'
  'D26' = '2.580.17'
  'D27' = '0.0₃0964'
  'D28' = '545.75'
  'D31' = '8.24'
  'D35' = '1.58'
  'D36' = '5.85'
  'D38' = '4.82'
  'D40' = '18.41'
  'D41' = '1.77'
  'D42' = '141.12'
  'D44' = '40.44'
  'D46' = '146.52'
  'D47' = '3.59'
  'D48' = '21.36'
  'D50' = '0.587'
  'D51' = '0.0936'
}
foreach ($addr in $dPrices.Keys) {
  $rng = $ws.Range($addr)
  $rng.NumberFormat = "@"
  $rng.Value = $dPrices[$addr]
  $rng.Style = "Normal"
}
